$wb = $excel.ActiveWorkbook

# --- Productdata sheet: raise AverageDemand for product 1 from 40 to 70 ---
$wsProductdata = $wb.Worksheets.Item("Productdata")
$wsProductdata.Range("G2").Value = 70

# Re-normalize the always-blank StandardDevDemands column so the
# round-trip through the engine doesn't turn its empty string cells
# into a stray "Name" value (keeps H2:H11 genuinely empty).
$wsProductdata.Range("H2:H11").ClearContents()

# --- ForecastedAverageDemand sheet: give positive demand to the last 3 periods ---
$wsAvgDemand = $wb.Worksheets.Item("ForecastedAverageDemand")
$wsAvgDemand.Range("B9").Value = 100
$wsAvgDemand.Range("B10").Value = 100
$wsAvgDemand.Range("B11").Value = 100

# --- ForcastedStandardDeviation sheet: matching std-dev values for those periods ---
$wsStdDev = $wb.Worksheets.Item("ForcastedStandardDeviation")
$wsStdDev.Range("B9").Value = 10.23775
$wsStdDev.Range("B10").Value = 11.713975
$wsStdDev.Range("B11").Value = 13.0425775
